# Append the new resale-numbers row for 2023-06-13 18:23:24 (row 43)
# to the CityResaleNum sheet, extending the used range to A1:T43.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$row = 43

# Columns A (Date) and D (Week) hold values that look numeric/date-like
# ("2023-06-13", "24"); a leading apostrophe forces them to stay text,
# matching the existing rows in the sheet (e.g. row 42).
$ws.Cells.Item($row, 1).Value = "'2023-06-13"
$ws.Cells.Item($row, 2).Value = "18:23:24"
$ws.Cells.Item($row, 3).Value = "Tuesday"
$ws.Cells.Item($row, 4).Value = "'24"
$ws.Cells.Item($row, 5).Value = 121115
$ws.Cells.Item($row, 6).Value = 135051
$ws.Cells.Item($row, 7).Value = 161213
$ws.Cells.Item($row, 8).Value = 132440
$ws.Cells.Item($row, 9).Value = 176360
$ws.Cells.Item($row, 10).Value = 113893
$ws.Cells.Item($row, 11).Value = 202320
$ws.Cells.Item($row, 12).Value = 223145
$ws.Cells.Item($row, 13).Value = 173959
$ws.Cells.Item($row, 14).Value = 101513
$ws.Cells.Item($row, 15).Value = 38845
$ws.Cells.Item($row, 16).Value = 34059
$ws.Cells.Item($row, 17).Value = 51432
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36804
$ws.Cells.Item($row, 20).Value = -1
